$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2 and 3 to the new date (2022-10-04 = 44838)
$ws.Range("D2").Value = 44838
$ws.Range("D3").Value = 44838

# Append the previous week's data as new rows 4 and 5 (2022-09-28 = 44832)
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 44832
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 100112039
$ws.Range("G4").Value = "Ciboulette"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 1200
$ws.Range("L4").Value = 1300
$ws.Range("M4").Value = 1250
$ws.Range("N4").Value = "`$/docena de atados"
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 417
$ws.Range("Q4").Value = 3
$ws.Range("R4").Value = "Hortaliza"

$ws.Range("D4").NumberFormat = $ws.Range("D2").NumberFormat

$ws.Range("A5").Value = 7
$ws.Range("B5").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C5").Value = "Ñuble"
$ws.Range("D5").Value = 44832
$ws.Range("E5").Value = 16
$ws.Range("F5").Value = 100112039
$ws.Range("G5").Value = "Ciboulette"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 1000
$ws.Range("M5").Value = 1000
$ws.Range("N5").Value = "`$/docena de atados"
$ws.Range("O5").Value = "Región Metropolitana"
$ws.Range("P5").Value = 333
$ws.Range("Q5").Value = 3
$ws.Range("R5").Value = "Hortaliza"

$ws.Range("D5").NumberFormat = $ws.Range("D3").NumberFormat
